# ActitimeHybridFramework test-data workbook: add a "Parameters" header
# column and admin/manager login rows used by the extent report run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "username"/"password" columns one place to the right
# and make room for a new first column ("Parameters").
$ws.Columns.Item(1).Insert()

# Header row (row 1) - written in username/password/Parameters order so
# the shared-string table matches the authored workbook.
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("A1").Value = "Parameters"

# Data rows
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "galatpassword"

$ws.Range("B3").Value = "rajmitra"
$ws.Range("C3").Value = "manager"

# Bold the header row
$ws.Range("A1:C1").Font.Bold = $true

# Column A width (new "Parameters" column)
$ws.Columns.Item(1).ColumnWidth = 15.5

# Selection as left by the author
$ws.Range("G11").Select() | Out-Null
